# Add season-record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1="Wins", AE1="Losses", AF1="Ties" using the same bold
# header style as the rest of row 1 (e.g. AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record for every player row (2-49): Seattle Mariners' 2013 record
# was 71 wins, 91 losses, 0 ties.
$wins = 71
$losses = 91
$ties = 0

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
